# tg-1178 missing views rename column names CLEANUP, ADD documentation
#
# Blad1 had a leftover scratch/test row ("Ikbestaniet" / "bla" / "bal") that
# was never meant to ship, plus the BGT_PND_pand view's attribute-documentation
# cell (F7) still carried the wrong column list copied from the "type" views
# instead of documenting its own identifying column. Clean both up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# CLEANUP: drop the bogus placeholder row 8 ("Ikbestaniet"/"bla"/"bal")
$ws.Rows.Item(8).Delete()

# ADD documentation: BGT_PND_pand's VIEW ATTRIBUTEN should document the
# identifying column, not the bgt_type/plus_type attribute list it
# mistakenly shared with the other "gebouw" views.
$ws.Range("F7").Value = "identificatie_lokaalid"

# Re-apply the autofilter over the (now verified) header/data range so the
# workbook's filter-database bookkeeping is refreshed/current.
$filterRange = $ws.Range("A1:J4")
$ws.Names.Add("_xlnm._FilterDatabase_0_0_0_0", $filterRange)
$ws.Names.Add("_xlnm._FilterDatabase_0_0_0_0_0", $filterRange)
